# Generate Report for Handback
#
# For both locale sheets (zh-cn, de-de), row 8 (the b0cb1f99-... file) gets
# a "handback received" reconciliation pass: the Latest Target File /
# Latest Handback File / Latest Handback DateTime columns get filled in,
# an error is recorded because the handback is based on an older commit
# than the current "latest" handoff, and column P (Error Detail) is widened
# so the message is readable. A hyperlink is added on the new "Latest
# Target File" cell pointing at the same (latest) commit used by column A.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afa7a00f785a84d865926362e8b7588b1ee6851a/e2e/b0cb1f99-704b-4b56-9ea9-2431075ea55f.md"
$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d582f5581ae2cc737642dc598827c9120f4045d1/e2e/b0cb1f99-704b-4b56-9ea9-2431075ea55f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afa7a00f785a84d865926362e8b7588b1ee6851a/e2e/b0cb1f99-704b-4b56-9ea9-2431075ea55f.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I8").Value = "b0cb1f99-704b-4b56-9ea9-2431075ea55f.md"
$wsZh.Range("J8").Value = "b0cb1f99-704b-4b56-9ea9-2431075ea55f.bba51d15c5663390ecba0e9950c6f8acad5fdbd2.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-09-02 10:51:17"
$wsZh.Range("P8").Value = $errorMessage

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestUrl, "", "", "b0cb1f99-704b-4b56-9ea9-2431075ea55f.md")

$wsZh.Columns.Item(16).ColumnWidth = 39.14285714285714

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I8").Value = "b0cb1f99-704b-4b56-9ea9-2431075ea55f.md"
$wsDe.Range("J8").Value = "b0cb1f99-704b-4b56-9ea9-2431075ea55f.bba51d15c5663390ecba0e9950c6f8acad5fdbd2.de-de.xlf"
$wsDe.Range("K8").Value = "2016-09-02 10:51:24"
$wsDe.Range("P8").Value = $errorMessage

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestUrl, "", "", "b0cb1f99-704b-4b56-9ea9-2431075ea55f.md")

$wsDe.Columns.Item(16).ColumnWidth = 39.14285714285714
